# Insert a new weekly record as row 29 (pushing the existing rows 29-64 down to 30-65),
# matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 29; everything below shifts down one row
# and the used range/dimension grows from R64 to R65 automatically.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new weekly observation.
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(29, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(29, 4).Value = '2022-04-28'
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112009
$ws.Cells.Item(29, 7).Value = 'Acelga'
$ws.Cells.Item(29, 8).Value = 'Sin especificar'
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 1800
$ws.Cells.Item(29, 12).Value = 2000
$ws.Cells.Item(29, 13).Value = 1900
$ws.Cells.Item(29, 14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(29, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(29, 16).Value = 633
$ws.Cells.Item(29, 17).Value = 3
$ws.Cells.Item(29, 18).Value = 'Hortaliza'
